# Endosos para polizas blanket y complementaria
#
# Inserts two new columns (CodigoAgente, NUM_GRUPO) right before the
# existing "NroCuenta" column, fills in their header + data values, clears
# the old "FechaInicio" data value (keeping the cell/style), and updates the
# "Kms" data value. Finally moves the active selection to F8 to match the
# saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before column E (old NroCuenta column), shifting
# every column from E onward two places to the right (E->G, F->H, ... R->T).
$ws.Columns("E:F").Insert() | Out-Null

# New header cells.
$ws.Range("E1").Value = "CodigoAgente"
$ws.Range("F1").Value = "NUM_GRUPO"

# New data cells for row 2.
$ws.Range("E2").Value = 6965
$ws.Range("F2").Value = "Corporativos Directos"

# Old FechaInicio data (now at S2 after the column insert) is cleared out,
# leaving just the formatted, empty cell behind.
$ws.Range("S2").ClearContents() | Out-Null

# Old Kms data (now at T2 after the column insert) changes value.
$ws.Range("T2").Value = 2500

# Match the saved selection state.
$ws.Range("F8").Select() | Out-Null
